# Phase 3 report: fix the "ConcetID" typo in the header row to "ConceptID",
# and update the sheet view (zoom + selected cell) to match how the sheet
# was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fix the misspelled header in A1 ("ConcetID" -> "ConceptID").
$ws.Range("A1").Value = "ConceptID"

# Update the view: zoomed to 130% and D16 selected (new active cell),
# matching the saved sheetView state.
$excel.ActiveWindow.Zoom = 130
$ws.Range("D16").Select()
